# Apply the commit's changes to the "economics" workbook.
# Commit message: "Aenderungen zum grossen Teil eingebaut"
#  - adds a new "prChange_pellet" price-change-factor row to gen_economics
#  - tweaks the MIP time_limit in further_parameters from 200 to 100
#  - leftover cursor/selection moves from the editing session on a couple
#    of other sheets
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) gen_economics: insert a new row 7 for the pellet price-change factor
#    (pushes the former rows 7-10 down to 8-11)
# ---------------------------------------------------------------------
$wsGen = $wb.Worksheets.Item("gen_economics")
$wsGen.Activate()

$wsGen.Rows.Item(7).Insert()

$wsGen.Cells.Item(7, 1).Value = "prChange_pellet"
$wsGen.Cells.Item(7, 2).Value = 1.03
$wsGen.Cells.Item(7, 3).Value = "-"
$wsGen.Cells.Item(7, 4).Value = "Price change factors per year for pellets"

# Column D was widened to fit the new, longer description text.
$wsGen.Columns.Item(4).ColumnWidth = 79.67

$wsGen.Range("B7").Select() | Out-Null

# ---------------------------------------------------------------------
# 2) gas_economics: selection left on the (whole) row 3
# ---------------------------------------------------------------------
$wsGas = $wb.Worksheets.Item("gas_economics")
$wsGas.Activate()
$wsGas.Rows.Item(3).Select() | Out-Null

# ---------------------------------------------------------------------
# 3) el_economics: selection left on the (whole) row 3
# ---------------------------------------------------------------------
$wsEl = $wb.Worksheets.Item("el_economics")
$wsEl.Activate()
$wsEl.Rows.Item(3).Select() | Out-Null

# ---------------------------------------------------------------------
# 4) further_parameters: time_limit changed from 200 to 100, selection D14
#    (this sheet ends up the active / tab-selected one, matching the file)
# ---------------------------------------------------------------------
$wsFurther = $wb.Worksheets.Item("further_parameters")
$wsFurther.Cells.Item(3, 2).Value = 100
$wsFurther.Activate()
$wsFurther.Range("D14").Select() | Out-Null
